$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

# --- New row 8: SearchClient test data -------------------------------
$ws.Range("A8").Value = "SearchClient"
$ws.Range("B8").Value = "AT"
$ws.Range("C8").Value = "Bohar"
$ws.Range("D8").Value = "Password@12"
$ws.Range("E8").Value = "Bohar"

# Reuse the existing column formatting (border style already used by
# rows 2-7) for the new row's B/C/D cells.
$ws.Range("B2").Copy()
$ws.Range("B8").PasteSpecial(-4122)

$ws.Range("C2").Copy()
$ws.Range("C8").PasteSpecial(-4122)

# A8/E8 get a left+right thin-border-only style (new style in the diff).
$rngA8 = $ws.Range("A8")
$rngA8.Borders.Item(7).LineStyle = 1
$rngA8.Borders.Item(10).LineStyle = 1

$ws.Range("A8").Copy()
$ws.Range("E8").PasteSpecial(-4122)

# D8 mirrors D2/D6/D7: a mailto hyperlink to the password value. Add the
# hyperlink first, then paste D2's formatting over it so D8 ends up with
# exactly the same style as the other hyperlinked password cells.
$ws.Hyperlinks.Add($ws.Range("D8"), "mailto:Password@12")

$ws.Range("D2").Copy()
$ws.Range("D8").PasteSpecial(-4122)

# --- Sheet view: move the active selection like in the diff ----------
[void]$ws.Range("K22").Select()
